# works_published_brill.xlsx — "Updating collaboration.R with 2022 data"
#
# The underlying OpenAlex/Brill export got re-sorted: rows 3 & 4 swap places
# (the "Mummy Labels" work moves above the "Consciousness Is Quantum State
# Reduction" work) and rows 7 & 8 swap places (the "Art, Allegory..." review
# moves above the "Doctrinal Engagements..." paper). Every column of each
# row-pair trades places verbatim (including which columns are populated at
# all, e.g. the abstract in column E), so this is a full row-content swap,
# not a per-field edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out how many columns are actually in use (dimension is A1:AL8 here)
# so the swap covers every populated column without hard-coding the sheet
# width.
$lastCol = $ws.UsedRange.Columns.Count()

# Columns whose text values look like plain numbers/dates ("2023-01-04",
# "158", "1-3", ...). Excel's COM layer auto-coerces such literals (dates in
# particular) into real date/number values when assigned via .Value, which
# would silently change the stored cell type away from the original text
# cells. Forcing a "Text" number format on the destination cell before the
# assignment keeps them as plain text, matching the source data.
$textCols = @(6, 15, 16, 17, 18)   # F=publication_date, O,P,Q,R = pages/volume/issue

function Swap-Rows([int]$r1, [int]$r2, [int]$lastCol, $textCols) {
    # Snapshot every cell in both rows first so the write-back in the loop
    # below can't clobber a value we still need to read.
    $vals1 = @{}
    $vals2 = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $vals1[$c] = $ws.Cells.Item($r1, $c).Value()
        $vals2[$c] = $ws.Cells.Item($r2, $c).Value()
    }

    for ($c = 1; $c -le $lastCol; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)

        if ($textCols -contains $c) {
            $cell1.NumberFormat = "@"
            $cell2.NumberFormat = "@"
        }

        $cell1.Value = $vals2[$c]
        $cell2.Value = $vals1[$c]
    }
}

Swap-Rows 3 4 $lastCol $textCols
Swap-Rows 7 8 $lastCol $textCols
